$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.535.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "'1.839.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.74%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'226.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'32.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.75%  "

$ws.Range("D9").Value = "'0.295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.60%  "

$ws.Range("D10").Value = "'0.0719"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.93%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").Value = "'2.104.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("D13").Value = "'1.837.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.69%  "

$ws.Range("D14").Value = "'11.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").Value = "'0.651"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.94%  "

$ws.Range("D16").Value = "'34.541.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "'4.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.07%  "

$ws.Range("D18").Value = "'69.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'253.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").Value = "'0.0₃0805"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.67%  "

$ws.Range("D21").Value = "'11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.24%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "'4.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.21%  "

$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").Value = "'161.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.54%  "

$ws.Range("D26").Value = "'16.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").Value = "'7.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.15%  "

$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").Value = "'0.0537"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.50%  "

$ws.Range("D31").Value = "'3.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("B32").Value = "Swop.fi"
$ws.Range("C32").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D32").Value = "'519.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +896.09%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "

$ws.Range("D34").Value = "'3.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("D35").Value = "'1.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.84%  "

$ws.Range("D36").Value = "'1.460.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.655"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.28%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("E39").Value = "  +4.59%  "

$ws.Range("D40").Value = "'0.984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.32%  "

$ws.Range("D41").Value = "'83.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").Value = "'2.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("E44").Value = "  +5.83%  "

$ws.Range("E45").Value = "  +7.47%  "

$ws.Range("D46").Value = "'2.000.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.96%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("D49").Value = "'12.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.31%  "

$ws.Range("D50").Value = "'106.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.61%  "

$ws.Range("E51").Value = "  +0.21%  "
